# "fix input team listing"
#
# Adds a second team-roster sheet ("Team1") after the existing "Team0"
# sheet, fills in two more names that were missing from Team0's "sde"
# column (C8/C9), and lays out Team1 with the same pandas-export shape
# (header row in B1:E1, a 0-based index column A, and the first new
# name going into B2).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Team0")

# --- Team0: two names were missing from column C (sde) ---
$ws1.Range("C8").Value = "('331_John Ho', 'sde2')"
$ws1.Range("C9").Value = "('336_Mrs. Claudia Thomas', 'sde2')"

# --- add the new "Team1" sheet right after "Team0" ---
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "Team1"

# Match Team0's cell formatting (bold/bordered/centered header style and
# the matching style on column A) without disturbing the shared style
# table - copy just the formats, not the values.
$ws1.Range("A2:A21").Copy()
$ws2.Range("A2:A21").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("B1:E1").Copy()
$ws2.Range("B1:E1").PasteSpecial(-4122)    # xlPasteFormats

# Header row - reuse the same category labels (sdn / sde / nfit / ne)
$ws2.Range("B1").Value = $ws1.Range("B1").Value2
$ws2.Range("C1").Value = $ws1.Range("C1").Value2
$ws2.Range("D1").Value = $ws1.Range("D1").Value2
$ws2.Range("E1").Value = $ws1.Range("E1").Value2

# Index column A: 0..19 down rows 2..21
for ($i = 0; $i -le 19; $i++) {
    $r = $i + 2
    $ws2.Range("A$r").Value = $i
}

# First roster entry for Team1 goes under "sdn" (column B)
$ws2.Range("B2").Value = "('371_Michelle Kirby', 'sdn1')"

# Keep "Team0" as the active/selected sheet
$ws1.Activate()
